# Generate Report for Archive
# - Update the "Status" value from "Ready for handoff" to "In Translation"
#   on the Overview sheet (columns E & F) and on each per-language sheet
#   (column C on "zh-cn" and "de-de").
# - Re-fit the width of the Status columns that held the old, longer text
#   so they match the shorter new value.

$wb = $excel.ActiveWorkbook

# --- Overview sheet (columns E and F hold the Status for each language) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value2 = "In Translation"
$wsOverview.Range("F2").Value2 = "In Translation"
$wsOverview.Range("E1:F1").EntireColumn.ColumnWidth = 12.5

# --- zh-cn sheet (column C holds Status) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value2 = "In Translation"
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = 12.5

# --- de-de sheet (column C holds Status) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value2 = "In Translation"
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = 12.5
